$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 0.99994318200069288
$ws.Range("E2").Value = 0.99994318200069288

# Row 3
$ws.Range("C3").Value = $true
$ws.Range("D3").Value = 0.002241240459353195
$ws.Range("E3").Value = 0.002241240459353195

# Row 4
$ws.Range("D4").Value = 0.008908871094672628
$ws.Range("E4").Value = 0.008908871094672628

# Row 5
$ws.Range("D5").Value = 0.00000005546239393289819
$ws.Range("E5").Value = 0.00000005546239393289819

# Row 6
$ws.Range("D6").Value = 0.001893309186933003
$ws.Range("E6").Value = 0.001893309186933003

# Row 7
$ws.Range("D7").Value = 0.008514477086576482
$ws.Range("E7").Value = 0.99148552291342351

# Row 8
$ws.Range("D8").Value = 0.99950971260042432
$ws.Range("E8").Value = 0.000490287399575684

# Row 9
$ws.Range("D9").Value = 0.99999999999999956
$ws.Range("E9").Value = 0.0000000000000004440892098500626

# Row 10
$ws.Range("D10").Value = 0.0000000003257435380001836
$ws.Range("E10").Value = 0.99999999967425646

# Row 11
$ws.Range("D11").Value = 0.99850411830763974
$ws.Range("E11").Value = 0.001495881692360257
$ws.Range("F11").Value = 3.64016342163085893
$ws.Range("G11").Value = 0.69999999999999996
